$wb = $excel.ActiveWorkbook

# --- ALC ---
$ALC = $wb.Worksheets.Item("ALC")
# row 12
$ALC.Range("H12").Value = 208.11111
$ALC.Range("I12").Value = 162.16667
$ALC.Range("J12").Value = 300
$ALC.Range("K12").Value = 162.16667
$ALC.Range("L12").Value = 300
$ALC.Range("M12").Value = 7.833329999999989
$ALC.Range("N12").Value = -640

# row 98
$ALC.Range("H98").Value = 3642.2856
$ALC.Range("I98").Value = 3299.2
$ALC.Range("K98").Value = 3299.2
$ALC.Range("M98").Value = -1801.2

# row 113
$ALC.Range("H113").Value = 3214.7
$ALC.Range("I113").Value = 2800
$ALC.Range("J113").Value = 3491.1667
$ALC.Range("K113").Value = 2800
$ALC.Range("L113").Value = 3491.1667
$ALC.Range("M113").Value = 454
$ALC.Range("N113").Value = -9999.1667

# row 122
$ALC.Range("H122").Value = 3642.2856
$ALC.Range("I122").Value = 3299.2
$ALC.Range("K122").Value = 9897.599999999999
$ALC.Range("M122").Value = -7447.599999999999

# row 138
$ALC.Range("H138").Value = 2251.3103
$ALC.Range("I138").Value = 1585.1765
$ALC.Range("J138").Value = 3195
$ALC.Range("K138").Value = 4755.529500000001
$ALC.Range("L138").Value = 9585
$ALC.Range("M138").Value = 384.4704999999994
$ALC.Range("N138").Value = -19865


# --- ARM ---
$ARM = $wb.Worksheets.Item("ARM")
# row 32
$ARM.Range("H32").Value = 18164.682
$ARM.Range("I32").Value = 18749.72
$ARM.Range("K32").Value = 18749.72
$ARM.Range("M32").Value = -18462.72

# row 61
$ARM.Range("H61").Value = 1588572.2
$ARM.Range("I61").Value = 2084536.5
$ARM.Range("K61").Value = 2084536.5
$ARM.Range("M61").Value = -2084324.5

# row 122
$ARM.Range("H122").Value = 2938.85
$ARM.Range("I122").Value = 2692.8823
$ARM.Range("K122").Value = 8078.646900000001
$ARM.Range("M122").Value = -5628.646900000001

# row 136
$ARM.Range("H136").Value = 1588572.2
$ARM.Range("I136").Value = 2084536.5
$ARM.Range("K136").Value = 6253609.5
$ARM.Range("M136").Value = -6251059.5


# --- BSM ---
$BSM = $wb.Worksheets.Item("BSM")
# row 21
$BSM.Range("H21").Value = 43591.5
$BSM.Range("J21").Value = 43591.5
$BSM.Range("L21").Value = 43591.5
$BSM.Range("N21").Value = -44063.5

# row 86
$BSM.Range("H86").Value = 1999.5
$BSM.Range("I86").Value = 1999
$BSM.Range("K86").Value = 1999
$BSM.Range("M86").Value = -876

# row 89
$BSM.Range("H89").Value = 1999.5
$BSM.Range("I89").Value = 1999
$BSM.Range("K89").Value = 9995
$BSM.Range("M89").Value = -4379

# row 105
$BSM.Range("H105").Value = 2033.238
$BSM.Range("I105").Value = 2031.7368
$BSM.Range("J105").Value = 2047.5
$BSM.Range("K105").Value = 2031.7368
$BSM.Range("L105").Value = 2047.5
$BSM.Range("M105").Value = -284.7367999999999
$BSM.Range("N105").Value = -5541.5

# row 122
$BSM.Range("H122").Value = 61250
$BSM.Range("J122").Value = 61250
$BSM.Range("L122").Value = 61250
$BSM.Range("N122").Value = -71050

# row 134
$BSM.Range("H134").Value = 1015714.06
$BSM.Range("I134").Value = 994106.1
$BSM.Range("J134").Value = 1145361.5
$BSM.Range("K134").Value = 2982318.3
$BSM.Range("L134").Value = 3436084.5
$BSM.Range("M134").Value = -2979783.3
$BSM.Range("N134").Value = -3441154.5


# --- CRP ---
$CRP = $wb.Worksheets.Item("CRP")
# row 7
$CRP.Range("H7").Value = 243.54546
$CRP.Range("I7").Value = 263.72726
$CRP.Range("J7").Value = 223.36363
$CRP.Range("K7").Value = 263.72726
$CRP.Range("L7").Value = 223.36363
$CRP.Range("M7").Value = -150.72726
$CRP.Range("N7").Value = -449.36363

# row 31
$CRP.Range("H31").Value = 28457.412
$CRP.Range("I31").Value = 11901.223
$CRP.Range("J31").Value = 47083.125
$CRP.Range("K31").Value = 11901.223
$CRP.Range("L31").Value = 47083.125
$CRP.Range("M31").Value = -11606.223
$CRP.Range("N31").Value = -47673.125

# row 34
$CRP.Range("H34").Value = 28457.412
$CRP.Range("I34").Value = 11901.223
$CRP.Range("J34").Value = 47083.125
$CRP.Range("K34").Value = 11901.223
$CRP.Range("L34").Value = 47083.125
$CRP.Range("M34").Value = -11699.223
$CRP.Range("N34").Value = -47487.125

# row 87
$CRP.Range("H87").Value = 118915
$CRP.Range("J87").Value = 118915
$CRP.Range("L87").Value = 118915
$CRP.Range("N87").Value = -121287

# row 90
$CRP.Range("H90").Value = 118915
$CRP.Range("J90").Value = 118915
$CRP.Range("L90").Value = 356745
$CRP.Range("N90").Value = -368601

# row 132
$CRP.Range("H132").Value = 27919090
$CRP.Range("I132").Value = 34484356
$CRP.Range("K132").Value = 103453068
$CRP.Range("M132").Value = -103450538


# --- CUL ---
$CUL = $wb.Worksheets.Item("CUL")
# row 4
$CUL.Range("H4").Value = 31966964
$CUL.Range("I4").Value = 38584164
$CUL.Range("K4").Value = 115752492
$CUL.Range("M4").Value = -115752380


# --- GSM ---
$GSM = $wb.Worksheets.Item("GSM")
# row 2
$GSM.Range("H2").Value = 719.1579
$GSM.Range("I2").Value = 966.5
$GSM.Range("K2").Value = 966.5
$GSM.Range("M2").Value = -853.5

# row 4
$GSM.Range("H4").Value = 1000
$GSM.Range("I4").Value = 0
$GSM.Range("J4").Value = 1000
$GSM.Range("K4").Value = 0
$GSM.Range("L4").Value = 1000
$GSM.Range("M4").ClearContents()
$GSM.Range("N4").Value = -1224

# row 113
$GSM.Range("H113").Value = 4492.971
$GSM.Range("I113").Value = 3808.2917
$GSM.Range("K113").Value = 3808.2917
$GSM.Range("M113").Value = -1638.2917

# row 124
$GSM.Range("H124").Value = 0
$GSM.Range("J124").Value = 0
$GSM.Range("L124").Value = 0
$GSM.Range("N124").ClearContents()

# row 126
$GSM.Range("H126").Value = 983306.8
$GSM.Range("J126").Value = 3440
$GSM.Range("L126").Value = 10320
$GSM.Range("N126").Value = -15260

# row 132
$GSM.Range("H132").Value = 37041210
$GSM.Range("I132").Value = 45457468
$GSM.Range("K132").Value = 136372404
$GSM.Range("M132").Value = -136369874


# --- LTW ---
$LTW = $wb.Worksheets.Item("LTW")
# row 68
$LTW.Range("H68").Value = 2499.5
$LTW.Range("I68").Value = 2499.25
$LTW.Range("J68").Value = 2500
$LTW.Range("K68").Value = 2499.25
$LTW.Range("L68").Value = 2500
$LTW.Range("M68").Value = -1750.25
$LTW.Range("N68").Value = -3998

# row 71
$LTW.Range("H71").Value = 2499.5
$LTW.Range("I71").Value = 2499.25
$LTW.Range("J71").Value = 2500
$LTW.Range("K71").Value = 12496.25
$LTW.Range("L71").Value = 12500
$LTW.Range("M71").Value = -8752.25
$LTW.Range("N71").Value = -19988

# row 132
$LTW.Range("H132").Value = 14137.125
$LTW.Range("I132").Value = 16331.105
$LTW.Range("K132").Value = 48993.315
$LTW.Range("M132").Value = -46463.315

# row 136
$LTW.Range("H136").Value = 4098.5
$LTW.Range("I136").Value = 2998.2856
$LTW.Range("J136").Value = 6665.6665
$LTW.Range("K136").Value = 8994.856800000001
$LTW.Range("L136").Value = 19996.9995
$LTW.Range("M136").Value = -6444.856800000001
$LTW.Range("N136").Value = -25096.9995


# --- WVR ---
$WVR = $wb.Worksheets.Item("WVR")
# row 58
$WVR.Range("H58").Value = 31299.8
$WVR.Range("I58").Value = 36874.75
$WVR.Range("K58").Value = 36874.75
$WVR.Range("M58").Value = -36566.75

# row 81
$WVR.Range("H81").Value = 4135.2144
$WVR.Range("I81").Value = 3689.9
$WVR.Range("J81").Value = 5248.5
$WVR.Range("K81").Value = 7379.8
$WVR.Range("L81").Value = 10497
$WVR.Range("M81").Value = -6318.8
$WVR.Range("N81").Value = -12619

# row 84
$WVR.Range("H84").Value = 4135.2144
$WVR.Range("I84").Value = 3689.9
$WVR.Range("J84").Value = 5248.5
$WVR.Range("K84").Value = 36899
$WVR.Range("L84").Value = 52485
$WVR.Range("M84").Value = -31595
$WVR.Range("N84").Value = -63093

# row 113
$WVR.Range("H113").Value = 1997.7858
$WVR.Range("I113").Value = 688.0714
$WVR.Range("K113").Value = 2064.2142
$WVR.Range("M113").Value = 105.7857999999997

# row 132
$WVR.Range("H132").Value = 21572124
$WVR.Range("I132").Value = 22919630
$WVR.Range("J132").Value = 12000
$WVR.Range("K132").Value = 68758890
$WVR.Range("L132").Value = 36000
$WVR.Range("M132").Value = -68756360
$WVR.Range("N132").Value = -41060

# row 133
$WVR.Range("H133").Value = 76666.336
$WVR.Range("J133").Value = 76666.336
$WVR.Range("L133").Value = 76666.336
$WVR.Range("N133").Value = -86786.336

# row 136
$WVR.Range("H136").Value = 14502776
$WVR.Range("I136").Value = 15161448
$WVR.Range("J136").Value = 12000
$WVR.Range("K136").Value = 45484344
$WVR.Range("L136").Value = 36000
$WVR.Range("M136").Value = -45481794
$WVR.Range("N136").Value = -41100

